$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 29.223446
$ws.Cells.Item(2, 8).Value = 87.670338
$ws.Cells.Item(2, 9).Value = 0.0169041244192178
$ws.Cells.Item(2, 10).Value = 0.0169041244192178
$ws.Cells.Item(2, 13).Value = 1.115279333333333
$ws.Cells.Item(2, 14).Value = 3.345838
$ws.Cells.Item(2, 15).Value = 0.08670767158519405
$ws.Cells.Item(2, 16).Value = 0.08670767158519403
$ws.Cells.Item(2, 17).Value = 32.59230537258266
$ws.Cells.Item(2, 18).Value = 293.330748353244
$ws.Cells.Item(2, 19).Value = 0.001465717268576796
$ws.Cells.Item(2, 20).Value = 0.001465717268576796

# Row 3
$ws.Cells.Item(3, 7).Value = 29.223446
$ws.Cells.Item(3, 8).Value = 87.670338
$ws.Cells.Item(3, 9).Value = 0.0169041244192178
$ws.Cells.Item(3, 10).Value = 0.0169041244192178
$ws.Cells.Item(3, 13).Value = 3.484068333333333
$ws.Cells.Item(3, 15).Value = 0.2708697667015328
$ws.Cells.Item(3, 16).Value = 0.2708697667015328
$ws.Cells.Item(3, 17).Value = 101.8164827994767
$ws.Cells.Item(3, 18).Value = 916.34834519529
$ws.Cells.Item(3, 19).Value = 0.00457881623772721
$ws.Cells.Item(3, 20).Value = 0.00457881623772721

# Row 4
$ws.Cells.Item(4, 7).Value = 29.223446
$ws.Cells.Item(4, 8).Value = 87.670338
$ws.Cells.Item(4, 9).Value = 0.0169041244192178
$ws.Cells.Item(4, 10).Value = 0.0169041244192178
$ws.Cells.Item(4, 13).Value = 8.022733000000001
$ws.Cells.Item(4, 14).Value = 24.068199
$ws.Cells.Item(4, 15).Value = 0.6237293899283516
$ws.Cells.Item(4, 16).Value = 0.6237293899283515
$ws.Cells.Item(4, 17).Value = 234.451904597918
$ws.Cells.Item(4, 18).Value = 2110.067141381262
$ws.Cells.Item(4, 19).Value = 0.01054359921127167
$ws.Cells.Item(4, 20).Value = 0.01054359921127167

# Row 5
$ws.Cells.Item(5, 7).Value = 29.223446
$ws.Cells.Item(5, 8).Value = 87.670338
$ws.Cells.Item(5, 9).Value = 0.0169041244192178
$ws.Cells.Item(5, 10).Value = 0.0169041244192178
$ws.Cells.Item(5, 13).Value = 0.2404413333333334
$ws.Cells.Item(5, 14).Value = 0.7213240000000001
$ws.Cells.Item(5, 15).Value = 0.0186931717849216
$ws.Cells.Item(5, 16).Value = 0.0186931717849216
$ws.Cells.Item(5, 17).Value = 7.026524320834667
$ws.Cells.Item(5, 18).Value = 63.23871888751201
$ws.Cells.Item(5, 19).Value = 0.0003159917016421265
$ws.Cells.Item(5, 20).Value = 0.0003159917016421264

# Row 6
$ws.Cells.Item(6, 9).Value = 0.9471112884046843
$ws.Cells.Item(6, 10).Value = 0.9471112884046842
$ws.Cells.Item(6, 13).Value = 1.115279333333333
$ws.Cells.Item(6, 14).Value = 3.345838
$ws.Cells.Item(6, 15).Value = 0.08670767158519405
$ws.Cells.Item(6, 16).Value = 0.08670767158519403
$ws.Cells.Item(6, 17).Value = 1826.095192390571
$ws.Cells.Item(6, 18).Value = 16434.85673151514
$ws.Cells.Item(6, 19).Value = 0.08212181454962336
$ws.Cells.Item(6, 20).Value = 0.08212181454962335

# Row 7
$ws.Cells.Item(7, 9).Value = 0.9471112884046843
$ws.Cells.Item(7, 10).Value = 0.9471112884046842
$ws.Cells.Item(7, 13).Value = 3.484068333333333
$ws.Cells.Item(7, 15).Value = 0.2708697667015328
$ws.Cells.Item(7, 16).Value = 0.2708697667015328
$ws.Cells.Item(7, 17).Value = 5704.616093301794
$ws.Cells.Item(7, 19).Value = 0.256543813730565
$ws.Cells.Item(7, 20).Value = 0.256543813730565

# Row 8
$ws.Cells.Item(8, 9).Value = 0.9471112884046843
$ws.Cells.Item(8, 10).Value = 0.9471112884046842
$ws.Cells.Item(8, 13).Value = 8.022733000000001
$ws.Cells.Item(8, 14).Value = 24.068199
$ws.Cells.Item(8, 15).Value = 0.6237293899283516
$ws.Cells.Item(8, 16).Value = 0.6237293899283515
$ws.Cells.Item(8, 17).Value = 13135.96847289066
$ws.Cells.Item(8, 18).Value = 118223.716256016
$ws.Cells.Item(8, 19).Value = 0.5907411461109088
$ws.Cells.Item(8, 20).Value = 0.5907411461109086

# Row 9
$ws.Cells.Item(9, 9).Value = 0.9471112884046843
$ws.Cells.Item(9, 10).Value = 0.9471112884046842
$ws.Cells.Item(9, 13).Value = 0.2404413333333334
$ws.Cells.Item(9, 14).Value = 0.7213240000000001
$ws.Cells.Item(9, 15).Value = 0.0186931717849216
$ws.Cells.Item(9, 16).Value = 0.0186931717849216
$ws.Cells.Item(9, 17).Value = 393.6850165955245
$ws.Cells.Item(9, 18).Value = 3543.16514935972
$ws.Cells.Item(9, 19).Value = 0.01770451401358719
$ws.Cells.Item(9, 20).Value = 0.01770451401358718

# Row 10
$ws.Cells.Item(10, 7).Value = 37.39212666666667
$ws.Cells.Item(10, 8).Value = 112.17638
$ws.Cells.Item(10, 9).Value = 0.02162924801792661
$ws.Cells.Item(10, 10).Value = 0.0216292480179266
$ws.Cells.Item(10, 13).Value = 1.115279333333333
$ws.Cells.Item(10, 14).Value = 3.345838
$ws.Cells.Item(10, 15).Value = 0.08670767158519405
$ws.Cells.Item(10, 16).Value = 0.08670767158519403
$ws.Cells.Item(10, 17).Value = 41.70266610071555
$ws.Cells.Item(10, 18).Value = 375.32399490644
$ws.Cells.Item(10, 19).Value = 0.00187542173377309
$ws.Cells.Item(10, 20).Value = 0.001875421733773089

# Row 11
$ws.Cells.Item(11, 7).Value = 37.39212666666667
$ws.Cells.Item(11, 8).Value = 112.17638
$ws.Cells.Item(11, 9).Value = 0.02162924801792661
$ws.Cells.Item(11, 10).Value = 0.0216292480179266
$ws.Cells.Item(11, 13).Value = 3.484068333333333
$ws.Cells.Item(11, 15).Value = 0.2708697667015328
$ws.Cells.Item(11, 16).Value = 0.2708697667015328
$ws.Cells.Item(11, 17).Value = 130.2767244353222
$ws.Cells.Item(11, 18).Value = 1172.4905199179
$ws.Cells.Item(11, 19).Value = 0.005858709364545371
$ws.Cells.Item(11, 20).Value = 0.00585870936454537

# Row 12
$ws.Cells.Item(12, 7).Value = 37.39212666666667
$ws.Cells.Item(12, 8).Value = 112.17638
$ws.Cells.Item(12, 9).Value = 0.02162924801792661
$ws.Cells.Item(12, 10).Value = 0.0216292480179266
$ws.Cells.Item(12, 13).Value = 8.022733000000001
$ws.Cells.Item(12, 14).Value = 24.068199
$ws.Cells.Item(12, 15).Value = 0.6237293899283516
$ws.Cells.Item(12, 16).Value = 0.6237293899283515
$ws.Cells.Item(12, 17).Value = 299.9870485488467
$ws.Cells.Item(12, 18).Value = 2699.88343693962
$ws.Cells.Item(12, 19).Value = 0.01349079767083037
$ws.Cells.Item(12, 20).Value = 0.01349079767083037

# Row 13
$ws.Cells.Item(13, 7).Value = 37.39212666666667
$ws.Cells.Item(13, 8).Value = 112.17638
$ws.Cells.Item(13, 9).Value = 0.02162924801792661
$ws.Cells.Item(13, 10).Value = 0.0216292480179266
$ws.Cells.Item(13, 13).Value = 0.2404413333333334
$ws.Cells.Item(13, 14).Value = 0.7213240000000001
$ws.Cells.Item(13, 15).Value = 0.0186931717849216
$ws.Cells.Item(13, 16).Value = 0.0186931717849216
$ws.Cells.Item(13, 17).Value = 8.990612791902224
$ws.Cells.Item(13, 18).Value = 80.91551512712002
$ws.Cells.Item(13, 19).Value = 0.0004043192487777771
$ws.Cells.Item(13, 20).Value = 0.000404319248777777

# Row 14
$ws.Cells.Item(14, 7).Value = 24.817167
$ws.Cells.Item(14, 8).Value = 74.45150100000001
$ws.Cells.Item(14, 9).Value = 0.01435533915817136
$ws.Cells.Item(14, 10).Value = 0.01435533915817136
$ws.Cells.Item(14, 13).Value = 1.115279333333333
$ws.Cells.Item(14, 14).Value = 3.345838
$ws.Cells.Item(14, 15).Value = 0.08670767158519405
$ws.Cells.Item(14, 16).Value = 0.08670767158519403
$ws.Cells.Item(14, 17).Value = 27.678073466982
$ws.Cells.Item(14, 18).Value = 249.102661202838
$ws.Cells.Item(14, 19).Value = 0.001244718033220798
$ws.Cells.Item(14, 20).Value = 0.001244718033220798

# Row 15
$ws.Cells.Item(15, 7).Value = 24.817167
$ws.Cells.Item(15, 8).Value = 74.45150100000001
$ws.Cells.Item(15, 9).Value = 0.01435533915817136
$ws.Cells.Item(15, 10).Value = 0.01435533915817136
$ws.Cells.Item(15, 13).Value = 3.484068333333333
$ws.Cells.Item(15, 15).Value = 0.2708697667015328
$ws.Cells.Item(15, 16).Value = 0.2708697667015328
$ws.Cells.Item(15, 17).Value = 86.46470566774499
$ws.Cells.Item(15, 18).Value = 778.182351009705
$ws.Cells.Item(15, 19).Value = 0.003888427368695255
$ws.Cells.Item(15, 20).Value = 0.003888427368695254

# Row 16
$ws.Cells.Item(16, 7).Value = 24.817167
$ws.Cells.Item(16, 8).Value = 74.45150100000001
$ws.Cells.Item(16, 9).Value = 0.01435533915817136
$ws.Cells.Item(16, 10).Value = 0.01435533915817136
$ws.Cells.Item(16, 13).Value = 8.022733000000001
$ws.Cells.Item(16, 14).Value = 24.068199
$ws.Cells.Item(16, 15).Value = 0.6237293899283516
$ws.Cells.Item(16, 16).Value = 0.6237293899283515
$ws.Cells.Item(16, 17).Value = 199.101504657411
$ws.Cells.Item(16, 18).Value = 1791.913541916699
$ws.Cells.Item(16, 19).Value = 0.008953846935340799
$ws.Cells.Item(16, 20).Value = 0.008953846935340798

# Row 17
$ws.Cells.Item(17, 7).Value = 24.817167
$ws.Cells.Item(17, 8).Value = 74.45150100000001
$ws.Cells.Item(17, 9).Value = 0.01435533915817136
$ws.Cells.Item(17, 10).Value = 0.01435533915817136
$ws.Cells.Item(17, 13).Value = 0.2404413333333334
$ws.Cells.Item(17, 14).Value = 0.7213240000000001
$ws.Cells.Item(17, 15).Value = 0.0186931717849216
$ws.Cells.Item(17, 16).Value = 0.0186931717849216
$ws.Cells.Item(17, 17).Value = 5.967072723036001
$ws.Cells.Item(17, 18).Value = 53.70365450732401
$ws.Cells.Item(17, 19).Value = 0.0002683468209145091
$ws.Cells.Item(17, 20).Value = 0.000268346820914509
